$d = $word.ActiveDocument

# Near the end of the document, in the "Requisitos" section, the body has:
#   ... "LOQ4073: Química Geral II (Requisito fraco)"        <- keep (last
#                                                                 requirement line)
#   (empty paragraph)                                         <- remove
#   "Ver no Jupiter Salvar em pdf Salvar em docx"             <- remove
#   "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll
#    and Github pages. Original theme under Creative Commons
#    Attribution"                                             <- remove
#   (empty paragraph)                                         <- keep
#   (page-break paragraph)                                    <- keep
#
# Locate the two footer text runs with Find, then grow the deletion range to
# also swallow the blank paragraph that precedes them and the paragraph mark
# that ends the "© 2020 ..." paragraph, so all three paragraphs disappear
# completely (not just their text).

$startRange = $d.Content.Duplicate
$found1 = $startRange.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$endRange = $d.Content.Duplicate
$found2 = $endRange.Find.Execute(
    "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found1 -or -not $found2) {
    throw "Could not locate the footer paragraphs to remove"
}

# The blank paragraph immediately before "Ver no Jupiter ..." is also removed.
$precedingBlank = $startRange.Paragraphs(1).Previous()

$delStart = $precedingBlank.Range.Start
$delEnd = $endRange.Paragraphs(1).Range.End

$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete()
